$d = $word.ActiveDocument
$para1 = $d.Paragraphs(1)
$end = $para1.Range.End - 1
$ip = $d.Range($end, $end)
$ip.InsertAfter(" for Initial Setup")
$newRange = $d.Range($end, $end + 19)
$newRange.Font.Bold = $true
$newRange.Font.Name = "Segoe UI"
$newRange.Font.NameBi = "Segoe UI"

$p1 = $d.Paragraphs(1)
$trueEnd = $p1.Range.End - 1
Write-Output ("trueEnd=$trueEnd")

# workaround for bookmark-add bug at paragraph End-1:
$tmpR = $d.Range($trueEnd, $trueEnd)
$tmpR.InsertAfter("X")
$bmR = $d.Range($trueEnd, $trueEnd)
$d.Bookmarks.Add("_GoBack", $bmR)
$delR = $d.Range($trueEnd, $trueEnd + 1)
Write-Output ("deleting=[" + $delR.Text + "]")
$delR.Delete()

$b = $d.Bookmarks("_GoBack")
Write-Output ("bm start=" + $b.Start + " end=" + $b.End)
Write-Output ("paragraph text=[" + $d.Paragraphs(1).Range.Text + "]")
